# Natmi following Dr Hou advice:
# the ligand/receptor summary now also reports the "ECs" sending-cluster
# (previously only FAPs/sCs were sent from), so the table grows from
# 2 sending clusters x 3 target clusters (6 data rows) to
# 3 sending clusters x 3 target clusters (9 data rows), and every
# edge-weight / specificity figure is recomputed against the larger
# ECs+FAPs+sCs population.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Timp3"
$ws.Range("C2").Value = "Kdr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 134.6741433333333
$ws.Range("H2").Value = 404.02243
$ws.Range("I2").Value = 0.4097716001282303
$ws.Range("J2").Value = 0.4097716001282303
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 162.98837
$ws.Range("N2").Value = 488.96511
$ws.Range("O2").Value = 0.9909539753179891
$ws.Range("P2").Value = 0.9909539753179891
$ws.Range("Q2").Value = 21950.31910304636
$ws.Range("R2").Value = 197552.8719274173
$ws.Range("S2").Value = 0.4060647961194833
$ws.Range("T2").Value = 0.4060647961194833

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Timp3"
$ws.Range("C3").Value = "Kdr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 134.6741433333333
$ws.Range("H3").Value = 404.02243
$ws.Range("I3").Value = 0.4097716001282303
$ws.Range("J3").Value = 0.4097716001282303
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6513563333333333
$ws.Range("N3").Value = 1.954069
$ws.Range("O3").Value = 0.003960185305646138
$ws.Range("P3").Value = 0.003960185305646138
$ws.Range("Q3").Value = 87.72085619640777
$ws.Range("R3").Value = 789.48770576767
$ws.Range("S3").Value = 0.001622771469498923
$ws.Range("T3").Value = 0.001622771469498923

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Timp3"
$ws.Range("C4").Value = "Kdr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 134.6741433333333
$ws.Range("H4").Value = 404.02243
$ws.Range("I4").Value = 0.4097716001282303
$ws.Range("J4").Value = 0.4097716001282303
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8364996666666666
$ws.Range("N4").Value = 2.509499
$ws.Range("O4").Value = 0.005085839376364744
$ws.Range("P4").Value = 0.005085839376364744
$ws.Range("Q4").Value = 112.6548760069522
$ws.Range("R4").Value = 1013.89388406257
$ws.Range("S4").Value = 0.002084032539248142
$ws.Range("T4").Value = 0.002084032539248142

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Timp3"
$ws.Range("C5").Value = "Kdr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 131.625918
$ws.Range("H5").Value = 394.877754
$ws.Range("I5").Value = 0.4004967969516487
$ws.Range("J5").Value = 0.4004967969516487
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 162.98837
$ws.Range("N5").Value = 488.96511
$ws.Range("O5").Value = 0.9909539753179891
$ws.Range("P5").Value = 0.9909539753179891
$ws.Range("Q5").Value = 21453.49382457366
$ws.Range("R5").Value = 193081.4444211629
$ws.Range("S5").Value = 0.3968738930413578
$ws.Range("T5").Value = 0.3968738930413578

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Timp3"
$ws.Range("C6").Value = "Kdr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 131.625918
$ws.Range("H6").Value = 394.877754
$ws.Range("I6").Value = 0.4004967969516487
$ws.Range("J6").Value = 0.4004967969516487
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6513563333333333
$ws.Range("N6").Value = 1.954069
$ws.Range("O6").Value = 0.003960185305646138
$ws.Range("P6").Value = 0.003960185305646138
$ws.Range("Q6").Value = 85.73537532011399
$ws.Range("R6").Value = 771.6183778810259
$ws.Range("S6").Value = 0.001586041530246264
$ws.Range("T6").Value = 0.001586041530246264

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Timp3"
$ws.Range("C7").Value = "Kdr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 131.625918
$ws.Range("H7").Value = 394.877754
$ws.Range("I7").Value = 0.4004967969516487
$ws.Range("J7").Value = 0.4004967969516487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8364996666666666
$ws.Range("N7").Value = 2.509499
$ws.Range("O7").Value = 0.005085839376364744
$ws.Range("P7").Value = 0.005085839376364744
$ws.Range("Q7").Value = 110.105036531694
$ws.Range("R7").Value = 990.9453287852459
$ws.Range("S7").Value = 0.002036862380044651
$ws.Range("T7").Value = 0.002036862380044651

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Timp3"
$ws.Range("C8").Value = "Kdr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 62.35654466666666
$ws.Range("H8").Value = 187.069634
$ws.Range("I8").Value = 0.189731602920121
$ws.Range("J8").Value = 0.189731602920121
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 162.98837
$ws.Range("N8").Value = 488.96511
$ws.Range("O8").Value = 0.9909539753179891
$ws.Range("P8").Value = 0.9909539753179891
$ws.Range("Q8").Value = 10163.39157405219
$ws.Range("R8").Value = 91470.52416646972
$ws.Range("S8").Value = 0.1880152861571481
$ws.Range("T8").Value = 0.1880152861571481

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Timp3"
$ws.Range("C9").Value = "Kdr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 62.35654466666666
$ws.Range("H9").Value = 187.069634
$ws.Range("I9").Value = 0.189731602920121
$ws.Range("J9").Value = 0.189731602920121
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6513563333333333
$ws.Range("N9").Value = 1.954069
$ws.Range("O9").Value = 0.003960185305646138
$ws.Range("P9").Value = 0.003960185305646138
$ws.Range("Q9").Value = 40.61633029341621
$ws.Range("R9").Value = 365.546972640746
$ws.Range("S9").Value = 0.0007513723059009513
$ws.Range("T9").Value = 0.0007513723059009513

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Timp3"
$ws.Range("C10").Value = "Kdr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 62.35654466666666
$ws.Range("H10").Value = 187.069634
$ws.Range("I10").Value = 0.189731602920121
$ws.Range("J10").Value = 0.189731602920121
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8364996666666666
$ws.Range("N10").Value = 2.509499
$ws.Range("O10").Value = 0.005085839376364744
$ws.Range("P10").Value = 0.005085839376364744
$ws.Range("Q10").Value = 52.16122882815177
$ws.Range("R10").Value = 469.4510594533659
$ws.Range("S10").Value = 0.0009649444570719516
$ws.Range("T10").Value = 0.0009649444570719516
